$d = $word.ActiveDocument

$d.Content.Find.Execute("30÷7=4, 2", $true, $false, $false, $false, $false, $true, 1, $false, "28÷8=3, 4", 2) | Out-Null
$d.Content.Find.Execute("42÷7=6, 0", $true, $false, $false, $false, $false, $true, 1, $false, "64÷6=10, 4", 2) | Out-Null
$d.Content.Find.Execute("45÷9=5, 0", $true, $false, $false, $false, $false, $true, 1, $false, "77÷8=9, 5", 2) | Out-Null
$d.Content.Find.Execute("84÷3=28, 0", $true, $false, $false, $false, $false, $true, 1, $false, "94÷2=47, 0", 2) | Out-Null
$d.Content.Find.Execute("82÷8=10, 2", $true, $false, $false, $false, $false, $true, 1, $false, "63÷4=15, 3", 2) | Out-Null
$d.Content.Find.Execute("35÷4=8, 3", $true, $false, $false, $false, $false, $true, 1, $false, "55÷4=13, 3", 2) | Out-Null
$d.Content.Find.Execute("94÷6=15, 4", $true, $false, $false, $false, $false, $true, 1, $false, "96÷6=16, 0", 2) | Out-Null
$d.Content.Find.Execute("25÷3=8, 1", $true, $false, $false, $false, $false, $true, 1, $false, "33÷2=16, 1", 2) | Out-Null
$d.Content.Find.Execute("36÷8=4, 4", $true, $false, $false, $false, $false, $true, 1, $false, "16÷3=5, 1", 2) | Out-Null
$d.Content.Find.Execute("19÷4=4, 3", $true, $false, $false, $false, $false, $true, 1, $false, "86÷6=14, 2", 2) | Out-Null
$d.Content.Find.Execute("65÷6=10, 5", $true, $false, $false, $false, $false, $true, 1, $false, "19÷5=3, 4", 2) | Out-Null
$d.Content.Find.Execute("21÷2=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "85÷8=10, 5", 2) | Out-Null
$d.Content.Find.Execute("66÷4=16, 2", $true, $false, $false, $false, $false, $true, 1, $false, "92÷9=10, 2", 2) | Out-Null
$d.Content.Find.Execute("98÷5=19, 3", $true, $false, $false, $false, $false, $true, 1, $false, "46÷4=11, 2", 2) | Out-Null
$d.Content.Find.Execute("76÷6=12, 4", $true, $false, $false, $false, $false, $true, 1, $false, "20÷8=2, 4", 2) | Out-Null
$d.Content.Find.Execute("81÷9=9, 0", $true, $false, $false, $false, $false, $true, 1, $false, "52÷8=6, 4", 2) | Out-Null
$d.Content.Find.Execute("24÷8=3, 0", $true, $false, $false, $false, $false, $true, 1, $false, "19÷4=4, 3", 2) | Out-Null
$d.Content.Find.Execute("65÷4=16, 1", $true, $false, $false, $false, $false, $true, 1, $false, "34÷5=6, 4", 2) | Out-Null
$d.Content.Find.Execute("13÷9=1, 4", $true, $false, $false, $false, $false, $true, 1, $false, "36÷6=6, 0", 2) | Out-Null
$d.Content.Find.Execute("43÷8=5, 3", $true, $false, $false, $false, $false, $true, 1, $false, "30÷8=3, 6", 2) | Out-Null
$d.Content.Find.Execute("28÷3=9, 1", $true, $false, $false, $false, $false, $true, 1, $false, "49÷2=24, 1", 2) | Out-Null
$d.Content.Find.Execute("74÷3=24, 2", $true, $false, $false, $false, $false, $true, 1, $false, "71÷9=7, 8", 2) | Out-Null
$d.Content.Find.Execute("57÷5=11, 2", $true, $false, $false, $false, $false, $true, 1, $false, "91÷4=22, 3", 2) | Out-Null
$d.Content.Find.Execute("41÷7=5, 6", $true, $false, $false, $false, $false, $true, 1, $false, "47÷5=9, 2", 2) | Out-Null
$d.Content.Find.Execute("13÷5=2, 3", $true, $false, $false, $false, $false, $true, 1, $false, "17÷9=1, 8", 2) | Out-Null
